$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: small in-place fixups (string re-index only; text unchanged or category correction) ----
$ws.Range("B1").Value = "التاريخ "
$ws.Range("D22").Value = "تغطية اعلامية"
$ws.Range("I22").Value = 2

# ---- Step 2: rewrite rows 23-29 values (row 23 previously held what is now row 28 data; rows 24-29 are new/restructured entries) ----
# Row 23
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 45668
$ws.Range("C23").Value = "الضيافة"
$ws.Range("D23").Value = "احتفالية معا نستطيع"
$ws.Range("E23").Value = 40
$ws.Range("F23").Value = 95000
$ws.Range("G23").Value = 20011928328
$ws.Range("H23").Value = "اشعار"
$ws.Range("I23").Value = 1
$ws.Range("O23").Value = "المجلس القومي للتدريب"

# Row 24
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 45668
$ws.Range("C24").Value = "اطباق"
$ws.Range("D24").Value = "احتفالية معا نستطيع"
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = 10000
$ws.Range("G24").Value = 20011892418
$ws.Range("H24").Value = "اشعار"
$ws.Range("I24").Value = 1

# Row 25
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 45668
$ws.Range("C25").Value = "ترحيل"
$ws.Range("D25").Value = "ترحيل الاحتفال"
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 7000
$ws.Range("G25").Value = 20011917072
$ws.Range("H25").Value = "اشعار"
$ws.Range("I25").Value = 3

# Row 26
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 45668
$ws.Range("C26").Value = "تزين القاعة"
$ws.Range("D26").Value = "احتفالية معا نستطيع"
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 150000
$ws.Range("G26").Value = "20011907880 - 200119155061"
$ws.Range("H26").Value = "اشعار"
$ws.Range("I26").Value = 1

# Row 27
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 45668
$ws.Range("C27").Value = "البقالة"
$ws.Range("D27").Value = "احتفالية معا نستطيع"
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 144000
$ws.Range("I27").Value = 1

# Row 28
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 45652
$ws.Range("C28").Value = "المجلس القومي للتدريب"
$ws.Range("D28").Value = "تراخيص المركز"
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 160000
$ws.Range("G28").Value = 20022963334
$ws.Range("H28").Value = "اشعار"
$ws.Range("I28").Value = 4

# Row 29
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 45678
$ws.Range("C29").Value = "تجديد ترخيص ادارة التدريب الموحد"
$ws.Range("D29").Value = "تراخيص المركز"
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 120000
$ws.Range("G29").Value = 20025165140
$ws.Range("H29").Value = "اشعار"
$ws.Range("I29").Value = 4

# ---- Step 3: formulas ----
$ws.Range("J23").Formula = "=F23/E23"
$ws.Range("F30").Formula = "=SUM(F2:F29)"

# ---- Step 4: style fixups (copy formats from representative cells that already carry the right cellXf) ----
$targets_2 = @("G23", "H23", "C24", "H24", "H25", "C26", "H26", "C27", "E27", "F27", "G27", "H27", "C28", "E28", "F28", "G28", "H28", "C29", "E29", "F29", "G29", "H29")
foreach ($addr in $targets_2) {
    $ws.Range("C2").Copy()
    $ws.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

$targets_7 = @("D24", "D26", "A27", "D27", "A28", "D28", "A29", "D29")
foreach ($addr in $targets_7) {
    $ws.Range("D2").Copy()
    $ws.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

$targets_13 = @("I24", "I25", "I26", "I27", "I28", "I29")
foreach ($addr in $targets_13) {
    $ws.Range("I2").Copy()
    $ws.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

$targets_9 = @("B25", "B26", "B27", "B28", "B29")
foreach ($addr in $targets_9) {
    $ws.Range("B2").Copy()
    $ws.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# ---- Step 5: selection (matches final author view) ----
$ws.Range("I30").Select()

$excel.CutCopyMode = $false
